# Add 5 more "fake records" for Gender = "f" (rows 7-11), mirroring the
# existing Adidas / "m" rows 2-6.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$brand = "Adidas"
$gender = "f"

# Row 7  <- mirrors row 2
$ws.Range("A7").Value = $brand
$ws.Range("B7").Value = $gender
$ws.Range("C7").Value = 35.5
$ws.Range("D7").Value = 5.5
$ws.Range("E7").Value = 4
$ws.Range("F7").Value = 37
$ws.Range("F7").NumberFormat = "# ?/?"

# Row 8  <- mirrors row 3
$ws.Range("A8").Value = $brand
$ws.Range("B8").Value = $gender
$ws.Range("C8").Value = 36
$ws.Range("D8").Value = 6
$ws.Range("E8").Value = 4.5
$ws.Range("F8").Formula = '="37 1/3"'
$ws.Range("F8").NumberFormat = "# ?/?"

# Row 9  <- mirrors row 4
$ws.Range("A9").Value = $brand
$ws.Range("B9").Value = $gender
$ws.Range("C9").Value = 36.5
$ws.Range("D9").Value = 6.5
$ws.Range("E9").Value = 5
$ws.Range("F9").Value = 38

# Row 10 <- mirrors row 5
$ws.Range("A10").Value = $brand
$ws.Range("B10").Value = $gender
$ws.Range("C10").Value = 37
$ws.Range("D10").Value = 7
$ws.Range("E10").Value = 5.5
$ws.Range("F10").Formula = '="38 1/3"'

# Row 11 <- mirrors row 6
$ws.Range("A11").Value = $brand
$ws.Range("B11").Value = $gender
$ws.Range("C11").Value = 37.5
$ws.Range("D11").Value = 7.5
$ws.Range("E11").Value = 6
$ws.Range("F11").Value = 39

# Update selection to match the authored state
$ws.Range("B11").Select()
